$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.894424333333333
$ws.Range("H2").Value = 5.683273
$ws.Range("I2").Value = 0.6002819911800915
$ws.Range("J2").Value = 0.6002819911800916
$ws.Range("O2").Value = 0.7806153008439968
$ws.Range("P2").Value = 0.7806153008439968
$ws.Range("Q2").Value = 0.09052822414222221
$ws.Range("R2").Value = 0.8147540172799999
$ws.Range("S2").Value = 0.4685893071362806
$ws.Range("T2").Value = 0.4685893071362807

$ws.Range("G3").Value = 1.894424333333333
$ws.Range("H3").Value = 5.683273
$ws.Range("I3").Value = 0.6002819911800915
$ws.Range("J3").Value = 0.6002819911800916
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01343
$ws.Range("N3").Value = 0.04029
$ws.Range("O3").Value = 0.2193846991560033
$ws.Range("P3").Value = 0.2193846991560033
$ws.Range("Q3").Value = 0.02544211879666666
$ws.Range("R3").Value = 0.22897906917
$ws.Range("S3").Value = 0.131692684043811
$ws.Range("T3").Value = 0.131692684043811

$ws.Range("I4").Value = 0.3198928944728968
$ws.Range("J4").Value = 0.3198928944728969
$ws.Range("O4").Value = 0.7806153008439968
$ws.Range("P4").Value = 0.7806153008439968
$ws.Range("S4").Value = 0.2497132880568173
$ws.Range("T4").Value = 0.2497132880568173

$ws.Range("I5").Value = 0.3198928944728968
$ws.Range("J5").Value = 0.3198928944728969
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01343
$ws.Range("N5").Value = 0.04029
$ws.Range("O5").Value = 0.2193846991560033
$ws.Range("P5").Value = 0.2193846991560033
$ws.Range("Q5").Value = 0.01355821621
$ws.Range("R5").Value = 0.12202394589
$ws.Range("S5").Value = 0.07017960641607958
$ws.Range("T5").Value = 0.07017960641607959

$ws.Range("G6").Value = 0.1798433333333334
$ws.Range("H6").Value = 0.5395300000000001
$ws.Range("I6").Value = 0.05698655382231239
$ws.Range("J6").Value = 0.05698655382231241
$ws.Range("O6").Value = 0.7806153008439968
$ws.Range("P6").Value = 0.7806153008439968
$ws.Range("Q6").Value = 0.008594113422222223
$ws.Range("R6").Value = 0.0773470208
$ws.Range("S6").Value = 0.044484575856067
$ws.Range("T6").Value = 0.04448457585606701

$ws.Range("G7").Value = 0.1798433333333334
$ws.Range("H7").Value = 0.5395300000000001
$ws.Range("I7").Value = 0.05698655382231239
$ws.Range("J7").Value = 0.05698655382231241
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01343
$ws.Range("N7").Value = 0.04029
$ws.Range("O7").Value = 0.2193846991560033
$ws.Range("P7").Value = 0.2193846991560033
$ws.Range("Q7").Value = 0.002415295966666667
$ws.Range("R7").Value = 0.0217376637
$ws.Range("S7").Value = 0.01250197796624539
$ws.Range("T7").Value = 0.0125019779662454

$ws.Range("G8").Value = 0.07207599999999999
$ws.Range("H8").Value = 0.216228
$ws.Range("I8").Value = 0.0228385605246992
$ws.Range("J8").Value = 0.02283856052469921
$ws.Range("O8").Value = 0.7806153008439968
$ws.Range("P8").Value = 0.7806153008439968
$ws.Range("Q8").Value = 0.003444271786666666
$ws.Range("R8").Value = 0.03099844608
$ws.Range("S8").Value = 0.0178281297948319
$ws.Range("T8").Value = 0.0178281297948319

$ws.Range("G9").Value = 0.07207599999999999
$ws.Range("H9").Value = 0.216228
$ws.Range("I9").Value = 0.0228385605246992
$ws.Range("J9").Value = 0.02283856052469921
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01343
$ws.Range("N9").Value = 0.04029
$ws.Range("O9").Value = 0.2193846991560033
$ws.Range("P9").Value = 0.2193846991560033
$ws.Range("Q9").Value = 0.0009679806799999998
$ws.Range("R9").Value = 0.008711826119999998
$ws.Range("S9").Value = 0.005010430729867307
$ws.Range("T9").Value = 0.005010430729867309
